# Scheduled market-data refresh: overwrite price/profit columns (H-N) per leve row.
# Values sourced from latest Universalis snapshot for the Cerberus data centre.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 1851.125
$ws.Range("J17").Value = 1851.125
$ws.Range("L17").Value = 5553.375
$ws.Range("N17").Value = -5889.375

# Row 62
$ws.Range("H62").Value = 41671000
$ws.Range("I62").Value = 125002250
$ws.Range("K62").Value = 125002250
$ws.Range("M62").Value = -125001626

# Row 64
$ws.Range("H64").Value = 7488.4
$ws.Range("I64").Value = 7148.3335
$ws.Range("K64").Value = 7148.3335
$ws.Range("M64").Value = -6900.3335

# Row 65
$ws.Range("H65").Value = 41671000
$ws.Range("I65").Value = 125002250
$ws.Range("K65").Value = 625011250
$ws.Range("M65").Value = -625008130

# Row 67
$ws.Range("H67").Value = 7488.4
$ws.Range("I67").Value = 7148.3335
$ws.Range("K67").Value = 7148.3335
$ws.Range("M67").Value = -6290.3335

# Row 132
$ws.Range("H132").Value = 4660.516
$ws.Range("I132").Value = 4649.2
$ws.Range("K132").Value = 13947.6
$ws.Range("M132").Value = -11417.6

# Row 138
$ws.Range("H138").Value = 3494.712
$ws.Range("I138").Value = 4931.8823
$ws.Range("J138").Value = 2913
$ws.Range("K138").Value = 14795.6469
$ws.Range("L138").Value = 8739
$ws.Range("M138").Value = -9655.6469
$ws.Range("N138").Value = -19019

$ws = $wb.Worksheets.Item("ARM")
# Row 110
$ws.Range("H110").Value = 321
$ws.Range("I110").Value = 314.875
$ws.Range("K110").Value = 314.875
$ws.Range("M110").Value = 1730.125

# Row 132
$ws.Range("H132").Value = 2973.6875
$ws.Range("I132").Value = 2972.4546
$ws.Range("J132").Value = 2987.25
$ws.Range("K132").Value = 8917.363799999999
$ws.Range("L132").Value = 8961.75
$ws.Range("M132").Value = -6387.363799999999
$ws.Range("N132").Value = -14021.75

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 5083.2617
$ws.Range("I86").Value = 2799.5518
$ws.Range("J86").Value = 10177.692
$ws.Range("K86").Value = 2799.5518
$ws.Range("L86").Value = 10177.692
$ws.Range("M86").Value = -1676.5518
$ws.Range("N86").Value = -12423.692

# Row 89
$ws.Range("H89").Value = 5083.2617
$ws.Range("I89").Value = 2799.5518
$ws.Range("J89").Value = 10177.692
$ws.Range("K89").Value = 13997.759
$ws.Range("L89").Value = 50888.45999999999
$ws.Range("M89").Value = -8381.759000000002
$ws.Range("N89").Value = -62120.45999999999

# Row 105
$ws.Range("H105").Value = 2634.3
$ws.Range("J105").Value = 5291.5
$ws.Range("L105").Value = 5291.5
$ws.Range("N105").Value = -8785.5

# Row 134
$ws.Range("H134").Value = 10157.934
$ws.Range("I134").Value = 12045.556
$ws.Range("K134").Value = 36136.66800000001
$ws.Range("M134").Value = -33601.66800000001

$ws = $wb.Worksheets.Item("CRP")
# Row 2
$ws.Range("H2").Value = 6333.6665
$ws.Range("I2").Value = 7001
$ws.Range("J2").Value = 6000
$ws.Range("K2").Value = 7001
$ws.Range("L2").Value = 6000
$ws.Range("M2").Value = -6888
$ws.Range("N2").Value = -6226

# Row 31
$ws.Range("H31").Value = 2175.1738
$ws.Range("I31").Value = 1904.8334
$ws.Range("J31").Value = 2470.0908
$ws.Range("K31").Value = 1904.8334
$ws.Range("L31").Value = 2470.0908
$ws.Range("M31").Value = -1609.8334
$ws.Range("N31").Value = -3060.0908

# Row 34
$ws.Range("H34").Value = 2175.1738
$ws.Range("I34").Value = 1904.8334
$ws.Range("J34").Value = 2470.0908
$ws.Range("K34").Value = 1904.8334
$ws.Range("L34").Value = 2470.0908
$ws.Range("M34").Value = -1702.8334
$ws.Range("N34").Value = -2874.0908

# Row 59
$ws.Range("H59").Value = 99999
$ws.Range("J59").Value = 99999
$ws.Range("L59").Value = 99999
$ws.Range("N59").Value = -102289

$ws = $wb.Worksheets.Item("CUL")
# Row 37
$ws.Range("H37").Value = 129247.61
$ws.Range("J37").Value = 129247.61
$ws.Range("L37").Value = 387742.83
$ws.Range("N37").Value = -387966.83

# Row 81
$ws.Range("H81").Value = 28333.334
$ws.Range("I81").Value = 50000
$ws.Range("K81").Value = 150000
$ws.Range("M81").Value = -148877

# Row 84
$ws.Range("H84").Value = 28333.334
$ws.Range("I84").Value = 50000
$ws.Range("K84").Value = 450000
$ws.Range("M84").Value = -444384

# Row 131
$ws.Range("H131").Value = 22285404
$ws.Range("J131").Value = 25737674
$ws.Range("L131").Value = 77213022
$ws.Range("N131").Value = -77223102

# Row 141
$ws.Range("H141").Value = 31553.562
$ws.Range("I141").Value = 8671.4
$ws.Range("K141").Value = 26014.2
$ws.Range("M141").Value = -20834.2

$ws = $wb.Worksheets.Item("GSM")
# Row 69
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()

# Row 70
$ws.Range("H70").Value = 5644.9443
$ws.Range("I70").Value = 5207.1113
$ws.Range("J70").Value = 6082.778
$ws.Range("K70").Value = 5207.1113
$ws.Range("L70").Value = 6082.778
$ws.Range("M70").Value = -4937.1113
$ws.Range("N70").Value = -6622.778

# Row 72
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()

# Row 73
$ws.Range("H73").Value = 5644.9443
$ws.Range("I73").Value = 5207.1113
$ws.Range("J73").Value = 6082.778
$ws.Range("K73").Value = 5207.1113
$ws.Range("L73").Value = 6082.778
$ws.Range("M73").Value = -4271.1113
$ws.Range("N73").Value = -7954.778

# Row 80
$ws.Range("H80").Value = 3235
$ws.Range("I80").Value = 2626.7144
$ws.Range("K80").Value = 2626.7144
$ws.Range("M80").Value = -1628.7144

# Row 83
$ws.Range("H83").Value = 3235
$ws.Range("I83").Value = 2626.7144
$ws.Range("K83").Value = 13133.572
$ws.Range("M83").Value = -8141.572

# Row 139
$ws.Range("H139").Value = 88518
$ws.Range("J139").Value = 88518
$ws.Range("L139").Value = 88518
$ws.Range("N139").Value = -98798

$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 1834.625
$ws.Range("I68").Value = 1541.2727
$ws.Range("J68").Value = 2480
$ws.Range("K68").Value = 1541.2727
$ws.Range("L68").Value = 2480
$ws.Range("M68").Value = -792.2727
$ws.Range("N68").Value = -3978

# Row 71
$ws.Range("H71").Value = 1834.625
$ws.Range("I71").Value = 1541.2727
$ws.Range("J71").Value = 2480
$ws.Range("K71").Value = 7706.363499999999
$ws.Range("L71").Value = 12400
$ws.Range("M71").Value = -3962.363499999999
$ws.Range("N71").Value = -19888

$ws = $wb.Worksheets.Item("WVR")
# Row 51
$ws.Range("H51").Value = 31249.75
$ws.Range("I51").Value = 25000
$ws.Range("K51").Value = 25000
$ws.Range("M51").Value = -24490

# Row 93
$ws.Range("H93").Value = 50001
$ws.Range("J93").Value = 50001
$ws.Range("L93").Value = 50001
$ws.Range("N93").Value = -54993
